# Remove the trailing "<br>" from the LabelR values (column F) in Sheet 1.
# The data table's header is in row 1; data rows follow below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$labelCol = 6  # column F = LabelR

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $labelCol)
    $value = $cell.Value2
    if ($value -ne $null -and $value.EndsWith("<br>")) {
        $newValue = $value.Substring(0, $value.Length - 4)
        $cell.Value2 = $newValue
    }
}
